$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '37.752.89'
$c.Style = 'Normal'
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  +0.26%  '
$c.Style = 'Normal'
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '2.081.59'
$c.Style = 'Normal'
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  -0.79%  '
$c.Style = 'Normal'
$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  +0.12%  '
$c.Style = 'Normal'
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '234.41'
$c.Style = 'Normal'
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  -0.29%  '
$c.Style = 'Normal'
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '0.625'
$c.Style = 'Normal'
$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  -0.12%  '
$c.Style = 'Normal'
$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '58.70'
$c.Style = 'Normal'
$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  +0.71%  '
$c.Style = 'Normal'
$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  +0.12%  '
$c.Style = 'Normal'
$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  +0.20%  '
$c.Style = 'Normal'
$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.0786'
$c.Style = 'Normal'
$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  +0.53%  '
$c.Style = 'Normal'
$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  +2.74%  '
$c.Style = 'Normal'
$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '2.388.82'
$c.Style = 'Normal'
$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  -0.29%  '
$c.Style = 'Normal'
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '14.85'
$c.Style = 'Normal'
$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  +1.92%  '
$c.Style = 'Normal'
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '21.06'
$c.Style = 'Normal'
$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  -1.54%  '
$c.Style = 'Normal'
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  -2.24%  '
$c.Style = 'Normal'
$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  +1.29%  '
$c.Style = 'Normal'
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '2.047.32'
$c.Style = 'Normal'
$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  -1.74%  '
$c.Style = 'Normal'
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '37.697.64'
$c.Style = 'Normal'
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  +0.38%  '
$c.Style = 'Normal'
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '6.16'
$c.Style = 'Normal'
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  -0.73%  '
$c.Style = 'Normal'
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '71.26'
$c.Style = 'Normal'
$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  +2.02%  '
$c.Style = 'Normal'
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  +0.99%  '
$c.Style = 'Normal'
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '228.59'
$c.Style = 'Normal'
$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  +0.68%  '
$c.Style = 'Normal'
$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  -0.14%  '
$c.Style = 'Normal'
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  -1.41%  '
$c.Style = 'Normal'
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '169.25'
$c.Style = 'Normal'
$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  +0.01%  '
$c.Style = 'Normal'
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '0.138'
$c.Style = 'Normal'
$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  +3.31%  '
$c.Style = 'Normal'
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '8.99'
$c.Style = 'Normal'
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  +0.63%  '
$c.Style = 'Normal'
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  +0.92%  '
$c.Style = 'Normal'
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  -2.26%  '
$c.Style = 'Normal'
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '0.120'
$c.Style = 'Normal'
$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  +1.64%  '
$c.Style = 'Normal'
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  +0.72%  '
$c.Style = 'Normal'
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '0.0631'
$c.Style = 'Normal'
$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  +1.67%  '
$c.Style = 'Normal'
$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  +1.46%  '
$c.Style = 'Normal'
$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  -2.20%  '
$c.Style = 'Normal'
$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  +2.57%  '
$c.Style = 'Normal'
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '3.38'
$c.Style = 'Normal'
$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  -4.16%  '
$c.Style = 'Normal'
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.Style = 'Normal'
$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  -3.84%  '
$c.Style = 'Normal'
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '0.0981'
$c.Style = 'Normal'
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  +2.27%  '
$c.Style = 'Normal'
$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  +0.82%  '
$c.Style = 'Normal'
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '0.0215'
$c.Style = 'Normal'
$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  +1.12%  '
$c.Style = 'Normal'
$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  -2.91%  '
$c.Style = 'Normal'
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '1.460.32'
$c.Style = 'Normal'
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  -1.64%  '
$c.Style = 'Normal'
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '4.32'
$c.Style = 'Normal'
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  +4.70%  '
$c.Style = 'Normal'
$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.Style = 'Normal'
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '16.50'
$c.Style = 'Normal'
$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  +6.19%  '
$c.Style = 'Normal'
$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  +1.54%  '
$c.Style = 'Normal'
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '7.44'
$c.Style = 'Normal'
$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  +1.96%  '
$c.Style = 'Normal'
$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  +0.38%  '
$c.Style = 'Normal'
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '2.273.84'
$c.Style = 'Normal'
$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  -0.42%  '
$c.Style = 'Normal'
